$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '45.100.81'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.47%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.430.26'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.37%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.70'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.64%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '104.23'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +9.14%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.517'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.65%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.529'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +9.84%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.82'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.27%  '
$ws.Range("E11").Value = '  +2.06%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.123'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.48%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.60'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.99%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.96'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.88%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.812.59'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.64%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.422.36'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.16%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.834'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.32%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '44.973.16'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.08%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.42'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.71%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.36'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.51%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0919'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.42%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.78'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '243.71'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.48%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.30'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.45%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.51'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.90%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.50'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.87%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.19'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -7.48%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.57'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.34%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.68'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.13%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '48.99'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.12%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.128'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +16.64%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.67'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +11.70%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.23'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.04%  '
$ws.Range("E35").Value = '  +0.12%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0766'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.82%  '
$ws.Range("E37").Value = '  +3.99%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.50'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.43%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '127.57'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.30%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.88'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.08%  '
$ws.Range("E41").Value = '  +2.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.18'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.35%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.16'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.50%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0291'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.19%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.940.22'
$ws.Range("D45").Style = "Normal"
$ws.Range("E46").Value = '  -0.56%  '
$ws.Range("E47").Value = '  +8.11%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.23'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.05%  '
$ws.Range("E49").Value = '  +18.23%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '76.15'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.54%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '54.28'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.56%  '
